$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# Row 2
Set-TextValue "D2" "90.478.21"
Set-TextValue "E2" "  -0.93%  "

# Row 3
Set-TextValue "D3" "3.113.31"
Set-TextValue "E3" "  -0.64%  "

# Row 4
Set-TextValue "E4" "  +0.19%  "

# Row 5
Set-TextValue "D5" "233.31"
Set-TextValue "E5" "  +5.98%  "

# Row 6
Set-TextValue "D6" "627.24"
Set-TextValue "E6" "  +0.82%  "

# Row 7
Set-TextValue "E7" "  +9.84%  "

# Row 8
Set-TextValue "D8" "0.355"
Set-TextValue "E8" "  -5.43%  "

# Row 9
Set-TextValue "E9" "  +0.00%  "

# Row 10
Set-TextValue "D10" "3.112.26"
Set-TextValue "E10" "  -0.57%  "

# Row 11
Set-TextValue "D11" "0.715"
Set-TextValue "E11" "  -2.17%  "

# Row 12
Set-TextValue "E12" "  +3.08%  "

# Row 13
Set-TextValue "D13" "35.59"
Set-TextValue "E13" "  +2.61%  "

# Row 14
Set-TextValue "E14" "  -3.34%  "

# Row 15
Set-TextValue "D15" "5.57"
Set-TextValue "E15" "  +2.64%  "

# Row 16
Set-TextValue "D16" "90.229.39"
Set-TextValue "E16" "  -0.97%  "

# Row 17
Set-TextValue "D17" "3.686.14"
Set-TextValue "E17" "  -0.71%  "

# Row 18
Set-TextValue "D18" "3.111.19"
Set-TextValue "E18" "  -0.58%  "

# Row 19
Set-TextValue "E19" "  -3.13%  "

# Row 20
Set-TextValue "D20" "14.29"
Set-TextValue "E20" "  +0.73%  "

# Row 21
Set-TextValue "D21" "0.0000211"
Set-TextValue "E21" "  -6.60%  "

# Row 22
Set-TextValue "D22" "448.10"
Set-TextValue "E22" "  +3.52%  "

# Row 23
Set-TextValue "D23" "8.97"
Set-TextValue "E23" "  +2.74%  "

# Row 24
Set-TextValue "E24" "  +4.19%  "

# Row 25
Set-TextValue "D25" "5.88"
Set-TextValue "E25" "  -2.72%  "

# Row 26
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D26" "12.21"
Set-TextValue "E26" "  -1.05%  "

# Row 27
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D27" "85.03"
Set-TextValue "E27" "  +1.12%  "

# Row 28
Set-TextValue "D28" "3.274.72"
Set-TextValue "E28" "  -0.96%  "

# Row 29
Set-TextValue "D29" "0.996"
Set-TextValue "E29" "  -0.28%  "

# Row 30
Set-TextValue "D30" "9.20"
Set-TextValue "E30" "  +2.29%  "

# Row 31
Set-TextValue "D31" "0.158"
Set-TextValue "E31" "  -6.61%  "

# Row 32
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D32" "0.987"
Set-TextValue "E32" "  +11.72%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D33" "27.57"
Set-TextValue "E33" "  +17.43%  "

# Row 34
Set-TextValue "D34" "0.196"
Set-TextValue "E34" "  +29.13%  "

# Row 35
Set-TextValue "D35" "0.151"
Set-TextValue "E35" "  +5.91%  "

# Row 36
Set-TextValue "D36" "504.38"
Set-TextValue "E36" "  -4.71%  "

# Row 37
$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D37" "1.92"
Set-TextValue "E37" "  +3.45%  "

# Row 38
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D38" "3.63"
Set-TextValue "E38" "  -6.20%  "

# Row 39
Set-TextValue "D39" "6.96"
Set-TextValue "E39" "  -4.03%  "

# Row 40
Set-TextValue "E40" "  +0.05%  "

# Row 41
Set-TextValue "D41" "0.423"
Set-TextValue "E41" "  +11.38%  "

# Row 42
Set-TextValue "D42" "22.18"
Set-TextValue "E42" "  -0.56%  "

# Row 44
Set-TextValue "D44" "0.0846"
Set-TextValue "E44" "  +10.08%  "

# Row 45
Set-TextValue "E45" "  +29.78%  "

# Row 46
Set-TextValue "D46" "1.96"
Set-TextValue "E46" "  +1.44%  "

# Row 47
Set-TextValue "D47" "0.696"
Set-TextValue "E47" "  +11.58%  "

# Row 48
Set-TextValue "D48" "148.90"
Set-TextValue "E48" "  +3.03%  "

# Row 49
Set-TextValue "D49" "4.55"
Set-TextValue "E49" "  +9.08%  "

# Row 50
Set-TextValue "E50" "  +3.18%  "

# Row 51
Set-TextValue "D51" "44.62"
Set-TextValue "E51" "  +0.73%  "

Write-Host "Applied all changes"